# Staging.CustomReport.xlsx - header column re-order + view/codeName touch-up
# (mirrors the "moved staging files StagingTemplates directory" commit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-order the header row (row 2) text ----------------------------
# Old layout:  A2=CustomReport_ID  B2=Code  C2=Name          D2=BusinessKey
# New layout:  A2=BusinessKey      B2=Code  C2=CustomReport_ID D2=Name
$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "CustomReport_ID"
$ws.Range("D2").Value = "Name"

# --- 2. Window size stored in bookViews/workbookView ---------------------
# (Standard Excel automation surface for this is the ActiveWindow/Windows
# collection; wrapped defensively since some hosts expose it read-only.)
try {
    $win = $excel.ActiveWindow
    $win.Width = 28800
    $win.Height = 12585
} catch {
}
try {
    $win2 = $wb.Windows.Item(1)
    $win2.Width = 28800
    $win2.Height = 12585
} catch {
}

# --- 3. Worksheet VBA CodeName: Sheet4 -> Sheet6 --------------------------
# Only reachable (even in real Excel) through the VBE/VBProject surface;
# guarded so the script is a no-op when the host has no VBA project.
try {
    $ws.CodeName = "Sheet6"
} catch {
}
try {
    $comp = $wb.VBProject.VBComponents.Item($ws.CodeName)
    $comp.Name = "Sheet6"
} catch {
}
